$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 191, shifting existing rows 191:299 down to 192:300
$ws.Rows.Item(191).Insert()

# Populate the newly inserted row 191 with its data
$ws.Range("A191").Value = 4
$ws.Range("B191").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C191").Value = "Los Lagos"
$ws.Range("D191").Value = 44606
$ws.Range("E191").Value = 10
$ws.Range("F191").Value = 100114013
$ws.Range("G191").Value = "Zanahoria"
$ws.Range("H191").Value = "Sin especificar"
$ws.Range("I191").Value = "Primera"
$ws.Range("J191").Value = 150
$ws.Range("K191").Value = 12000
$ws.Range("L191").Value = 12000
$ws.Range("M191").Value = 12000
$ws.Range("N191").Value = "$/saco 20 kilos"
$ws.Range("O191").Value = "Chillán"
$ws.Range("P191").Value = 600
$ws.Range("Q191").Value = 20
$ws.Range("R191").Value = "Hortaliza"
